# Weekly update: insert a new "Apio" price observation for
# Terminal Hortofrutícola Agro Chillán right after the existing row 63,
# pushing the subsequent historical rows (old 64..104) down by one
# (new 65..105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64; Excel shifts rows 64..104 down to 65..105
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row with the new week's data point.
$ws.Cells.Item(64, 1).Value  = 7
$ws.Cells.Item(64, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(64, 3).Value  = "Ñuble"
$ws.Cells.Item(64, 4).Value  = 44438
$ws.Cells.Item(64, 5).Value  = 16
$ws.Cells.Item(64, 6).Value  = 100112017
$ws.Cells.Item(64, 7).Value  = "Apio"
$ws.Cells.Item(64, 8).Value  = "Americana (o)"
$ws.Cells.Item(64, 9).Value  = "Primera"
$ws.Cells.Item(64, 10).Value = 160
$ws.Cells.Item(64, 11).Value = 8000
$ws.Cells.Item(64, 12).Value = 9000
$ws.Cells.Item(64, 13).Value = 8500
$ws.Cells.Item(64, 14).Value = "`$/docena de matas"
$ws.Cells.Item(64, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(64, 16).Value = 1417
$ws.Cells.Item(64, 17).Value = 6
$ws.Cells.Item(64, 18).Value = "Hortaliza"
